# B6-PowerPoint.pptx edit:
#  1. Three tables (slides 14, 15, 16) switch from the bespoke
#     "Table_0" table style to the built-in "Medium Style 2 - Accent 1"
#     table style.
#  2. The deck's theme swaps from the "Integral" design (Red Violet
#     colour scheme) back to the default "Office Theme" colour scheme.

$p = $ppt.ActivePresentation

# ---- 1. Table styles -------------------------------------------------
$targetTableStyle = "{916EA1F1-D4A7-4BC2-A8A2-E3E5C7E99520}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($targetTableStyle, $true)
        }
    }
}

# ---- 2. Theme colours --------------------------------------------------
function ConvertTo-BgrInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - the stock "Office" theme
# palette that this deck originally shipped with before "Integral" was
# applied.
$officeThemeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-BgrInt $officeThemeColors[$i - 1]
}
